# Fruta / hortaliza, semanal
# Insert a new data row above current row 77 (pushing existing rows 77-81 down to 78-82),
# and populate the new row 77 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 77; this shifts rows 77:81 down to 78:82
$ws.Rows.Item(77).Insert()

# Fill in the new row 77 with its values
$ws.Range("A77").Value = 4
$ws.Range("B77").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C77").Value = "Los Lagos"
$ws.Range("D77").Value = [DateTime]::FromOADate(44585)
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = 100112026
$ws.Range("G77").Value = "Haba"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 70
$ws.Range("K77").Value = 22000
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = 22000
$ws.Range("N77").Value = "$/saco 25 kilos"
$ws.Range("O77").Value = "Región de La Araucanía"
$ws.Range("P77").Value = 880
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"
